$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "26.715.96"
$ws.Range("E2").Value = "  -0.02%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.635.28"
$ws.Range("E3").Value = "  -0.69%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.12%  "

# Row 5 - BNB
$ws.Range("D5").Value = "217.05"
$ws.Range("E5").Value = "  +0.51%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -1.02%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.19%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.81%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.0622"
$ws.Range("E9").Value = "  -0.72%  "

# Row 10 - Solana
$ws.Range("D10").Value = "18.99"
$ws.Range("E10").Value = "  -0.98%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  +0.19%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.863.48"
$ws.Range("E12").Value = "  -0.69%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.637.46"
$ws.Range("E13").Value = "  -0.26%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.17%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -1.57%  "

# Row 16 - Litecoin
$ws.Range("E16").Value = "  -1.48%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "26.702.90"
$ws.Range("E17").Value = "  -0.09%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -2.34%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.18%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "210.37"
$ws.Range("E20").Value = "  -3.73%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -0.81%  "

# Row 22 - now Toncoin (was Chainlink)
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").Value = "2.34"
$ws.Range("E22").Value = "  +2.07%  "

# Row 23 - now Chainlink (was Toncoin)
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "6.17"
$ws.Range("E23").Value = "  -1.60%  "

# Row 24 - Avalanche
$ws.Range("D24").Value = "9.24"
$ws.Range("E24").Value = "  -2.94%  "

# Row 25 - Monero
$ws.Range("D25").Value = "145.64"
$ws.Range("E25").Value = "  -0.20%  "

# Row 26 - BinanceUSD
$ws.Range("E26").Value = "  +0.07%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -2.18%  "

# Row 28 - Cosmos
$ws.Range("D28").Value = "7.04"
$ws.Range("E28").Value = "  -1.26%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -1.15%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.0503"
$ws.Range("E30").Value = "  -2.60%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.59%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  -0.41%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  -1.62%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.273.43"
$ws.Range("E34").Value = "  -0.57%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -1.65%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  +0.51%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  -2.11%  "

# Row 38 - ImmutableX
$ws.Range("D38").Value = "0.528"
$ws.Range("E38").Value = "  -1.57%  "

# Row 39 - ARBITRUM
$ws.Range("E39").Value = "  -2.32%  "

# Row 40 - PaxDollar
$ws.Range("E40").Value = "  +0.18%  "

# Row 41 - TrustWalletToken
$ws.Range("E41").Value = "  -1.60%  "

# Row 42 - MXToken
$ws.Range("E42").Value = "  -2.24%  "

# Row 43 - RocketPoolETH
$ws.Range("D43").Value = "1.774.20"
$ws.Range("E43").Value = "  -0.71%  "

# Row 44 - FraxShare
$ws.Range("E44").Value = "  -3.64%  "

# Row 45 - Aave
$ws.Range("D45").Value = "60.42"
$ws.Range("E45").Value = "  +0.97%  "

# Row 46 - Quant
$ws.Range("D46").Value = "91.08"
$ws.Range("E46").Value = "  -0.96%  "

# Row 47 - RenderToken
$ws.Range("E47").Value = "  -2.41%  "

# Row 48 - Cronos
$ws.Range("D48").Value = "0.0520"
$ws.Range("E48").Value = "  +0.79%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "7.52"
$ws.Range("E49").Value = "  -2.96%  "

# Row 50 - Algorand
$ws.Range("E50").Value = "  -0.79%  "

# Row 51 - now USDD (was Mantle)
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.01"
$ws.Range("E51").Value = "  +0.16%  "
